$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Valor Mora" total and "Cant. Periodos" count ------------------
$ws.Range("E11").Value = 512460
$ws.Range("F13").Value = 9

# --- Insert a new detail row for period 2509 -------------------------------
# The table of worker/period rows currently spans B16:J23 (periods 2501-2508).
# A new row is inserted right after the last data row (row 23) and it is
# filled with the same layout/format as that row, then the period value is
# changed to the new one (2509).
$ws.Rows("24:24").Insert()
$ws.Range("B23:J23").Copy($ws.Range("B24:J24"))
$ws.Range("E24").Value = "2509"

# Keep the "Observaciones"/signature block (now rows 29-30) looking right;
# Excel's row-insert already shifted it down from rows 28-29 automatically.

# --- Cosmetic: widen columns slightly (mirrors the bestFit recalculation
# that Excel performs once the new, slightly different, row is added) -------
$ws.Columns("B").ColumnWidth = 18.54
$ws.Columns("C").ColumnWidth = 16.73
$ws.Columns("D").ColumnWidth = 30.54
$ws.Columns("E").ColumnWidth = 13.54
$ws.Columns("F").ColumnWidth = 10.18
$ws.Columns("G").ColumnWidth = 14.36
$ws.Columns("H").ColumnWidth = 19.36
$ws.Columns("I").ColumnWidth = 18.09
$ws.Columns("J").ColumnWidth = 15
